$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AD2").Value = 5
$ws.Range("AE2").Value = 8
$ws.Range("AG2").Value = 8
$ws.Range("AM2").Value = 12
$ws.Range("AN2").Value = 8
$ws.Range("AV2").Value = 5
$ws.Range("AW2").Value = 11
$ws.Range("BA2").Value = 20
$ws.Range("BB2").Value = 19
$ws.Range("BF2").Value = "2012-04-21"

# Row 3
$ws.Range("AD3").Value = 1
$ws.Range("AF3").Value = 11
$ws.Range("AG3").Value = 11
$ws.Range("AH3").Value = 21
$ws.Range("AO3").Value = 25
$ws.Range("AU3").Value = 3
$ws.Range("AV3").Value = 18
$ws.Range("AX3").Value = 8
$ws.Range("BF3").Value = "2012-04-21"

# Row 4
$ws.Range("AD4").Value = 21
$ws.Range("AH4").Value = 25
$ws.Range("AO4").Value = 16
$ws.Range("AV4").Value = 12
$ws.Range("BF4").Value = "2012-04-21"

# Row 5
$ws.Range("D5").Value = 63
$ws.Range("E5").Value = 47
$ws.Range("G5").Value = 0.746
$ws.Range("I5").Value = 37.2
$ws.Range("J5").Value = 82.7
$ws.Range("L5").Value = 6.3
$ws.Range("N5").Value = 0.375
$ws.Range("O5").Value = 15.5
$ws.Range("P5").Value = 21.5
$ws.Range("R5").Value = 14.1
$ws.Range("S5").Value = 32.4
$ws.Range("U5").Value = 22.9
$ws.Range("AA5").Value = 18
$ws.Range("AC5").Value = 7.8
$ws.Range("AD5").Value = 5
$ws.Range("AH5").Value = 16
$ws.Range("AK5").Value = 13
$ws.Range("AM5").Value = 19
$ws.Range("AP5").Value = 20
$ws.Range("AU5").Value = 5
$ws.Range("BF5").Value = "2012-04-21"

# Row 6
$ws.Range("AD6").Value = 21
$ws.Range("AF6").Value = 24
$ws.Range("AH6").Value = 8
$ws.Range("AJ6").Value = 17
$ws.Range("AL6").Value = 13
$ws.Range("AN6").Value = 14
$ws.Range("AV6").Value = 25
$ws.Range("AX6").Value = 30
$ws.Range("BF6").Value = "2012-04-21"

# Row 7
$ws.Range("D7").Value = 64
$ws.Range("F7").Value = 28
$ws.Range("G7").Value = 0.5629999999999999
$ws.Range("I7").Value = 36.5
$ws.Range("K7").Value = 0.444
$ws.Range("M7").Value = 22.4
$ws.Range("O7").Value = 15.5
$ws.Range("Q7").Value = 0.769
$ws.Range("U7").Value = 21.1
$ws.Range("X7").Value = 5.2
$ws.Range("Y7").Value = 4.1
$ws.Range("Z7").Value = 18.8
$ws.Range("AA7").Value = 18.5
$ws.Range("AB7").Value = 96.09999999999999
$ws.Range("AC7").Value = 1.4
$ws.Range("AF7").Value = 12
$ws.Range("AG7").Value = 12
$ws.Range("AJ7").Value = 12
$ws.Range("AO7").Value = 26
$ws.Range("AQ7").Value = 12
$ws.Range("AU7").Value = 13
$ws.Range("AV7").Value = 8
$ws.Range("AX7").Value = 13
$ws.Range("BB7").Value = 18
$ws.Range("BF7").Value = "2012-04-21"

# Row 8
$ws.Range("D8").Value = 62
$ws.Range("E8").Value = 34
$ws.Range("G8").Value = 0.548
$ws.Range("I8").Value = 38.5
$ws.Range("J8").Value = 81.59999999999999
$ws.Range("K8").Value = 0.471
$ws.Range("M8").Value = 20.2
$ws.Range("N8").Value = 0.329
$ws.Range("O8").Value = 20
$ws.Range("P8").Value = 27.1
$ws.Range("Q8").Value = 0.735
$ws.Range("R8").Value = 11.2
$ws.Range("T8").Value = 43.1
$ws.Range("V8").Value = 15.5
$ws.Range("Z8").Value = 19.7
$ws.Range("AA8").Value = 22.3
$ws.Range("AB8").Value = 103.5
$ws.Range("AC8").Value = 1.9
$ws.Range("AD8").Value = 21
$ws.Range("AG8").Value = 13
$ws.Range("AN8").Value = 25
$ws.Range("AQ8").Value = 25
$ws.Range("AR8").Value = 16
$ws.Range("AT8").Value = 7
$ws.Range("AV8").Value = 27
$ws.Range("AW8").Value = 9
$ws.Range("AZ8").Value = 16
$ws.Range("BF8").Value = "2012-04-21"

# Row 9
$ws.Range("AD9").Value = 5
$ws.Range("AH9").Value = 9
$ws.Range("AN9").Value = 16
$ws.Range("AQ9").Value = 16
$ws.Range("AS9").Value = 29
$ws.Range("AU9").Value = 28
$ws.Range("BA9").Value = 17
$ws.Range("BF9").Value = "2012-04-21"

# Row 10
$ws.Range("D10").Value = 62
$ws.Range("F10").Value = 40
$ws.Range("G10").Value = 0.355
$ws.Range("J10").Value = 82.5
$ws.Range("L10").Value = 8.199999999999999
$ws.Range("M10").Value = 20.9
$ws.Range("N10").Value = 0.39
$ws.Range("O10").Value = 14.6
$ws.Range("Q10").Value = 0.771
$ws.Range("U10").Value = 22.5
$ws.Range("V10").Value = 13.8
$ws.Range("Z10").Value = 21.7
$ws.Range("AA10").Value = 16.6
$ws.Range("AD10").Value = 21
$ws.Range("AF10").Value = 22
$ws.Range("AJ10").Value = 9
$ws.Range("AK10").Value = 8
$ws.Range("AM10").Value = 8
$ws.Range("AU10").Value = 7
$ws.Range("AV10").Value = 4
$ws.Range("AX10").Value = 7
$ws.Range("BF10").Value = "2012-04-21"

# Row 11
$ws.Range("D11").Value = 63
$ws.Range("E11").Value = 32
$ws.Range("G11").Value = 0.508
$ws.Range("K11").Value = 0.45
$ws.Range("M11").Value = 19.9
$ws.Range("O11").Value = 15.7
$ws.Range("P11").Value = 20
$ws.Range("Q11").Value = 0.784
$ws.Range("R11").Value = 11.6
$ws.Range("S11").Value = 30.5
$ws.Range("T11").Value = 42.1
$ws.Range("W11").Value = 7.4
$ws.Range("AA11").Value = 18.6
$ws.Range("AC11").Value = 0.2
$ws.Range("AD11").Value = 5
$ws.Range("AE11").Value = 17
$ws.Range("AF11").Value = 18
$ws.Range("AG11").Value = 18
$ws.Range("AK11").Value = 14
$ws.Range("AM11").Value = 12
$ws.Range("AO11").Value = 22
$ws.Range("BA11").Value = 24
$ws.Range("BF11").Value = "2012-04-21"

# Row 12
$ws.Range("D12").Value = 63
$ws.Range("F12").Value = 22
$ws.Range("G12").Value = 0.651
$ws.Range("H12").Value = 48.4
$ws.Range("I12").Value = 35.6
$ws.Range("J12").Value = 81
$ws.Range("N12").Value = 0.37
$ws.Range("R12").Value = 12.4
$ws.Range("T12").Value = 43.7
$ws.Range("V12").Value = 14.1
$ws.Range("Y12").Value = 5.9
$ws.Range("Z12").Value = 21.8
$ws.Range("AB12").Value = 97.7
$ws.Range("AC12").Value = 3.5
$ws.Range("AD12").Value = 5
$ws.Range("AH12").Value = 9
$ws.Range("AJ12").Value = 18
$ws.Range("AU12").Value = 29
$ws.Range("AX12").Value = 6
$ws.Range("BF12").Value = "2012-04-21"

# Row 13
$ws.Range("AD13").Value = 5
$ws.Range("AH13").Value = 9
$ws.Range("AO13").Value = 23
$ws.Range("AP13").Value = 12
$ws.Range("AU13").Value = 14
$ws.Range("AZ13").Value = 26
$ws.Range("BB13").Value = 15
$ws.Range("BF13").Value = "2012-04-21"

# Row 14
$ws.Range("AD14").Value = 1
$ws.Range("AM14").Value = 18
$ws.Range("AQ14").Value = 15
$ws.Range("AU14").Value = 8
$ws.Range("BF14").Value = "2012-04-21"

# Row 15
$ws.Range("D15").Value = 63
$ws.Range("E15").Value = 38
$ws.Range("G15").Value = 0.603
$ws.Range("I15").Value = 36.7
$ws.Range("J15").Value = 82.2
$ws.Range("M15").Value = 12.8
$ws.Range("N15").Value = 0.33
$ws.Range("O15").Value = 17.3
$ws.Range("P15").Value = 22.8
$ws.Range("Q15").Value = 0.757
$ws.Range("R15").Value = 12.5
$ws.Range("S15").Value = 29.5
$ws.Range("V15").Value = 14.6
$ws.Range("W15").Value = 9.6
$ws.Range("Y15").Value = 5.6
$ws.Range("Z15").Value = 19.9
$ws.Range("AA15").Value = 19.4
$ws.Range("AC15").Value = 1.7
$ws.Range("AD15").Value = 5
$ws.Range("AE15").Value = 8
$ws.Range("AH15").Value = 16
$ws.Range("AJ15").Value = 11
$ws.Range("AN15").Value = 24
$ws.Range("AO15").Value = 10
$ws.Range("AQ15").Value = 17
$ws.Range("AR15").Value = 6
$ws.Range("AT15").Value = 17
$ws.Range("BA15").Value = 19
$ws.Range("BF15").Value = "2012-04-21"

# Row 16
$ws.Range("D16").Value = 62
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 0.726
$ws.Range("I16").Value = 37.5
$ws.Range("K16").Value = 0.474
$ws.Range("M16").Value = 15.4
$ws.Range("P16").Value = 24.6
$ws.Range("Q16").Value = 0.773
$ws.Range("R16").Value = 10.3
$ws.Range("S16").Value = 31.2
$ws.Range("V16").Value = 14.8
$ws.Range("W16").Value = 9
$ws.Range("X16").Value = 5.4
$ws.Range("AA16").Value = 20.6
$ws.Range("AB16").Value = 99.7
$ws.Range("AC16").Value = 7
$ws.Range("AD16").Value = 21
$ws.Range("AF16").Value = 3
$ws.Range("AO16").Value = 5
$ws.Range("AT16").Value = 21
$ws.Range("AV16").Value = 17
$ws.Range("BA16").Value = 8
$ws.Range("BB16").Value = 4
$ws.Range("BF16").Value = "2012-04-21"

# Row 17
$ws.Range("D17").Value = 62
$ws.Range("E17").Value = 29
$ws.Range("G17").Value = 0.468
$ws.Range("M17").Value = 19.5
$ws.Range("O17").Value = 16.5
$ws.Range("P17").Value = 21.1
$ws.Range("Q17").Value = 0.779
$ws.Range("S17").Value = 29.6
$ws.Range("T17").Value = 42.1
$ws.Range("AB17").Value = 99.59999999999999
$ws.Range("AC17").Value = 0.3
$ws.Range("AD17").Value = 21
$ws.Range("AI17").Value = 3
$ws.Range("AK17").Value = 18
$ws.Range("AL17").Value = 13
$ws.Range("AO17").Value = 17
$ws.Range("AP17").Value = 24
$ws.Range("AR17").Value = 5
$ws.Range("AS17").Value = 23
$ws.Range("AU17").Value = 2
$ws.Range("AV17").Value = 6
$ws.Range("AX17").Value = 17
$ws.Range("BA17").Value = 18
$ws.Range("BB17").Value = 5
$ws.Range("BF17").Value = "2012-04-21"

# Row 18
$ws.Range("AD18").Value = 1
$ws.Range("AH18").Value = 21
$ws.Range("AJ18").Value = 10
$ws.Range("AQ18").Value = 11
$ws.Range("BF18").Value = "2012-04-21"

# Row 19
$ws.Range("D19").Value = 63
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 0.349
$ws.Range("I19").Value = 34.4
$ws.Range("J19").Value = 80.59999999999999
$ws.Range("K19").Value = 0.426
$ws.Range("M19").Value = 22.8
$ws.Range("O19").Value = 17
$ws.Range("P19").Value = 21.9
$ws.Range("Q19").Value = 0.776
$ws.Range("R19").Value = 12
$ws.Range("S19").Value = 28.4
$ws.Range("AC19").Value = -5.3
$ws.Range("AD19").Value = 5
$ws.Range("AF19").Value = 24
$ws.Range("AG19").Value = 24
$ws.Range("AI19").Value = 28
$ws.Range("AN19").Value = 17
$ws.Range("AY19").Value = 17
$ws.Range("BB19").Value = 22
$ws.Range("BC19").Value = 25
$ws.Range("BF19").Value = "2012-04-21"

# Row 20
$ws.Range("AD20").Value = 5
$ws.Range("AV20").Value = 26
$ws.Range("BA20").Value = 25
$ws.Range("BF20").Value = "2012-04-21"

# Row 21
$ws.Range("AD21").Value = 5
$ws.Range("AE21").Value = 14
$ws.Range("AG21").Value = 14
$ws.Range("AI21").Value = 21
$ws.Range("BB21").Value = 14
$ws.Range("BF21").Value = "2012-04-21"

# Row 22
$ws.Range("AD22").Value = 5
$ws.Range("BF22").Value = "2012-04-21"

# Row 23
$ws.Range("D23").Value = 62
$ws.Range("F23").Value = 26
$ws.Range("G23").Value = 0.581
$ws.Range("H23").Value = 48.3
$ws.Range("I23").Value = 34.5
$ws.Range("J23").Value = 78
$ws.Range("L23").Value = 10.1
$ws.Range("M23").Value = 26.8
$ws.Range("N23").Value = 0.378
$ws.Range("O23").Value = 15.2
$ws.Range("P23").Value = 23.1
$ws.Range("Q23").Value = 0.658
$ws.Range("S23").Value = 31.5
$ws.Range("V23").Value = 15
$ws.Range("X23").Value = 4.1
$ws.Range("AA23").Value = 19.9
$ws.Range("AB23").Value = 94.5
$ws.Range("AC23").Value = 1.5
$ws.Range("AD23").Value = 21
$ws.Range("AG23").Value = 10
$ws.Range("AH23").Value = 14
$ws.Range("AP23").Value = 11
$ws.Range("AV23").Value = 19
$ws.Range("BA23").Value = 13
$ws.Range("BF23").Value = "2012-04-21"

# Row 24
$ws.Range("D24").Value = 62
$ws.Range("E24").Value = 32
$ws.Range("G24").Value = 0.516
$ws.Range("J24").Value = 83.2
$ws.Range("M24").Value = 14.5
$ws.Range("O24").Value = 13.3
$ws.Range("P24").Value = 18
$ws.Range("Q24").Value = 0.741
$ws.Range("R24").Value = 10.4
$ws.Range("S24").Value = 32.7
$ws.Range("U24").Value = 21.9
$ws.Range("V24").Value = 11.2
$ws.Range("X24").Value = 5.1
$ws.Range("Z24").Value = 17.7
$ws.Range("AA24").Value = 16.1
$ws.Range("AB24").Value = 93.40000000000001
$ws.Range("AC24").Value = 4.5
$ws.Range("AD24").Value = 21
$ws.Range("AE24").Value = 17
$ws.Range("AG24").Value = 17
$ws.Range("AH24").Value = 25
$ws.Range("AN24").Value = 7
$ws.Range("AR24").Value = 23
$ws.Range("AT24").Value = 8
$ws.Range("AW24").Value = 12
$ws.Range("AX24").Value = 14
$ws.Range("BB24").Value = 24
$ws.Range("BF24").Value = "2012-04-21"

# Row 25
$ws.Range("D25").Value = 63
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 0.524
$ws.Range("I25").Value = 37.7
$ws.Range("J25").Value = 82.09999999999999
$ws.Range("L25").Value = 6.7
$ws.Range("M25").Value = 19.6
$ws.Range("N25").Value = 0.343
$ws.Range("O25").Value = 16.3
$ws.Range("P25").Value = 21.4
$ws.Range("Q25").Value = 0.759
$ws.Range("R25").Value = 10.8
$ws.Range("T25").Value = 41.6
$ws.Range("AB25").Value = 98.3
$ws.Range("AC25").Value = 0.2
$ws.Range("AD25").Value = 5
$ws.Range("AE25").Value = 14
$ws.Range("AF25").Value = 14
$ws.Range("AG25").Value = 14
$ws.Range("AJ25").Value = 13
$ws.Range("AL25").Value = 15
$ws.Range("AN25").Value = 18
$ws.Range("AP25").Value = 21
$ws.Range("AQ25").Value = 14
$ws.Range("AV25").Value = 7
$ws.Range("BF25").Value = "2012-04-21"

# Row 26
$ws.Range("D26").Value = 63
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 0.444
$ws.Range("I26").Value = 36.5
$ws.Range("J26").Value = 82
$ws.Range("K26").Value = 0.445
$ws.Range("M26").Value = 20.9
$ws.Range("N26").Value = 0.349
$ws.Range("O26").Value = 17.2
$ws.Range("P26").Value = 21.6
$ws.Range("R26").Value = 11.1
$ws.Range("T26").Value = 40.7
$ws.Range("U26").Value = 20.5
$ws.Range("X26").Value = 5
$ws.Range("Y26").Value = 4.6
$ws.Range("AA26").Value = 19.9
$ws.Range("AB26").Value = 97.5
$ws.Range("AC26").Value = 0
$ws.Range("AD26").Value = 5
$ws.Range("AH26").Value = 9
$ws.Range("AI26").Value = 16
$ws.Range("AK26").Value = 17
$ws.Range("AM26").Value = 9
$ws.Range("AN26").Value = 13
$ws.Range("AO26").Value = 11
$ws.Range("AR26").Value = 18
$ws.Range("AS26").Value = 22
$ws.Range("AW26").Value = 13
$ws.Range("AX26").Value = 16
$ws.Range("BA26").Value = 14
$ws.Range("BB26").Value = 13
$ws.Range("BF26").Value = "2012-04-21"

# Row 27
$ws.Range("AD27").Value = 5
$ws.Range("AH27").Value = 24
$ws.Range("AM27").Value = 11
$ws.Range("AO27").Value = 12
$ws.Range("AQ27").Value = 24
$ws.Range("AT27").Value = 9
$ws.Range("AZ27").Value = 15
$ws.Range("BF27").Value = "2012-04-21"

# Row 28
$ws.Range("AD28").Value = 21
$ws.Range("AH28").Value = 14
$ws.Range("AT28").Value = 11
$ws.Range("AU28").Value = 4
$ws.Range("AY28").Value = 16
$ws.Range("BF28").Value = "2012-04-21"

# Row 29
$ws.Range("AD29").Value = 5
$ws.Range("AF29").Value = 24
$ws.Range("AG29").Value = 24
$ws.Range("AH29").Value = 9
$ws.Range("AI29").Value = 27
$ws.Range("AQ29").Value = 10
$ws.Range("AR29").Value = 22
$ws.Range("AT29").Value = 20
$ws.Range("BF29").Value = "2012-04-21"

# Row 30
$ws.Range("D30").Value = 63
$ws.Range("E30").Value = 33
$ws.Range("G30").Value = 0.524
$ws.Range("I30").Value = 38.2
$ws.Range("J30").Value = 83.7
$ws.Range("K30").Value = 0.456
$ws.Range("M30").Value = 12.8
$ws.Range("N30").Value = 0.322
$ws.Range("Q30").Value = 0.753
$ws.Range("V30").Value = 14.3
$ws.Range("Z30").Value = 22
$ws.Range("AB30").Value = 99.40000000000001
$ws.Range("AC30").Value = 0.3
$ws.Range("AD30").Value = 5
$ws.Range("AI30").Value = 4
$ws.Range("AK30").Value = 9
$ws.Range("AN30").Value = 28
$ws.Range("AO30").Value = 6
$ws.Range("AV30").Value = 13
$ws.Range("AW30").Value = 8
$ws.Range("BA30").Value = 7
$ws.Range("BB30").Value = 6
$ws.Range("BF30").Value = "2012-04-21"

# Row 31
$ws.Range("D31").Value = 62
$ws.Range("E31").Value = 16
$ws.Range("G31").Value = 0.258
$ws.Range("M31").Value = 16.3
$ws.Range("N31").Value = 0.322
$ws.Range("O31").Value = 15.6
$ws.Range("P31").Value = 21.4
$ws.Range("R31").Value = 11.8
$ws.Range("S31").Value = 29.8
$ws.Range("X31").Value = 6.5
$ws.Range("Z31").Value = 21.2
$ws.Range("AB31").Value = 93.40000000000001
$ws.Range("AC31").Value = -6.3
$ws.Range("AD31").Value = 21
$ws.Range("AN31").Value = 27
$ws.Range("AO31").Value = 24
$ws.Range("AP31").Value = 22
$ws.Range("AZ31").Value = 25
$ws.Range("BB31").Value = 23
$ws.Range("BF31").Value = "2012-04-21"
